$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally held yearly rows for 2008..2020 (rows 2..14).
# The update drops the oldest two years (2008, 2009) and appends a new
# year (2021) at the end, so the data now spans 2010..2021 (rows 2..13).

# Remove the 2008 and 2009 rows; this shifts 2010..2020 up to rows 2..12
# and shrinks the used range accordingly.
$ws.Rows("2:3").Delete()

# Append the new 2021 row as row 13.
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 7177
$ws.Range("C13").Value = 4050
$ws.Range("D13").Value = 93891
$ws.Range("E13").Value = 10399
$ws.Range("F13").Value = 508.86
$ws.Range("G13").Value = 4221.69
$ws.Range("H13").Value = 101.88
$ws.Range("I13").Value = 1188.54
$ws.Range("J13").Value = 342916.5
$ws.Range("K13").Value = 11306
$ws.Range("L13").Value = 33.33295
$ws.Range("M13").Value = 359
$ws.Range("N13").Value = 4921
$ws.Range("O13").Value = 26821
$ws.Range("P13").Value = ""
$ws.Range("Q13").Value = 1327
$ws.Range("R13").Value = 15508
$ws.Range("S13").Value = 58
